# MAI_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer note (A10)
#    from 2021-04-05 to 2021-04-06
#  - refresh the Weight / Percent Change figures for each holding
#    (rows 2-7, columns D & E)
#
# The worksheet is protected (legacy password hash "D382"), so we must
# unprotect before writing and re-protect afterwards to leave the sheet
# in the same protected state it started in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("D382")

# --- Disclaimer date bump -------------------------------------------------
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."

# --- Holdings weight / percent-change refresh -----------------------------
$ws.Range("D2").Value = 0.4883894097314794
$ws.Range("E2").Value = 0.003514252245216731

$ws.Range("D3").Value = 0.3338675934144865
$ws.Range("E3").Value = -0.0007816316560820269

$ws.Range("D4").Value = 0.09391853108797726
$ws.Range("E4").Value = 0.00542360201982417

$ws.Range("D5").Value = 0.05464508297759559
$ws.Range("E5").Value = 0.001721960739295092

$ws.Range("D6").Value = 0.02917938278846124
$ws.Range("E6").Value = 0.00515962592712027

$ws.Range("E7").Value = 0.00220939022209965

$ws.Protect("D382")
